$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.070.33'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '2.006.88'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '257.48'
$ws.Range('E5').Value = '  +4.41%  '
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.96'
$ws.Range('E8').Value = '  -7.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.378'
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0767'
$ws.Range('E10').Value = '  -4.97%  '
$ws.Range('E11').Value = '  -2.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.28'
$ws.Range('E12').Value = '  -4.69%  '
$ws.Range('D13').Value = '2.307.72'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.35'
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.795'
$ws.Range('E15').Value = '  -6.31%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.122.34'
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.19'
$ws.Range('E17').Value = '  -4.42%  '
$ws.Range('D18').Value = '36.972.59'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.84'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').Value = '0.0₃0831'
$ws.Range('E20').Value = '  -3.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '233.71'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.08'
$ws.Range('E22').Value = '  -2.56%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.75'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.93'
$ws.Range('E27').Value = '  -4.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.51'
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.35'
$ws.Range('E29').Value = '  -3.47%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.124'
$ws.Range('E30').Value = '  -9.62%  '
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.57'
$ws.Range('E32').Value = '  -3.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0640'
$ws.Range('E33').Value = '  -4.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.41'
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.35'
$ws.Range('E35').Value = '  -6.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.46'
$ws.Range('E36').Value = '  -3.92%  '
$ws.Range('E37').Value = '  +0.81%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.55'
$ws.Range('E39').Value = '  +4.32%  '
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').Value = '1.438.71'
$ws.Range('E42').Value = '  +4.46%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0920'
$ws.Range('E43').Value = '  -5.32%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0210'
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '89.03'
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.53'
$ws.Range('E46').Value = '  -7.21%  '
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.93'
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.94'
$ws.Range('E49').Value = '  -6.46%  '
$ws.Range('D50').Value = '2.199.91'
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.93'
$ws.Range('E51').Value = '  -8.04%  '
